# Mjolnir Group Presentation - fix factual error in the data-sources
# bullet on the "Introduction and assumptions" slide:
#   "The Meteoritical Society" -> "The Meteorological Society"

$p = $ppt.ActivePresentation

# sldId 257 (cId 485914459) is the 2nd slide in the deck.
$s = $p.Slides.Item(2)

# Shape id=3, "Content Placeholder 2" -> Shapes.Item(5) on this slide.
$sh = $s.Shapes.Item(5)
$tr = $sh.TextFrame.TextRange

$oldStr = "comprehensive data set of Meteorite Landings from The Meteoritical Society"
$newStr = "comprehensive data set of Meteorite Landings from The Meteorological Society"

$full = $tr.Text
$idx0 = $full.IndexOf($oldStr)
if ($idx0 -ge 0) {
    # Characters() is 1-based; replace the whole matched run text in one
    # shot so the existing run/formatting is preserved (only the <a:t/>
    # content changes).
    $run = $tr.Characters($idx0 + 1, $oldStr.Length)
    $run.Text = $newStr
} else {
    # Fallback: the exact phrase wasn't found verbatim (e.g. differing
    # whitespace) - just fix the one differing word directly.
    $word0 = $full.IndexOf("Meteoritical")
    if ($word0 -ge 0) {
        $word = $tr.Characters($word0 + 1, "Meteoritical".Length)
        $word.Text = "Meteorological"
    }
}
